$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 286
$ws.Range("F12").Value = 119
$ws.Range("F13").Value = 2498
$ws.Range("F14").Value = 62
$ws.Range("F15").Value = 29
$ws.Range("F17").Value = 17
$ws.Range("F20").Value = 593
$ws.Range("F21").Value = 177
$ws.Range("F22").Value = 92
$ws.Range("F25").Value = 2095
$ws.Range("F26").Value = 4191
$ws.Range("F30").Value = 1229
$ws.Range("F31").Value = 243
$ws.Range("F32").Value = 2130
$ws.Range("F34").Value = 472
$ws.Range("F36").Value = 126
$ws.Range("F38").Value = 438
$ws.Range("F39").Value = 725
$ws.Range("F42").Value = 7

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 47

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 286
$ws.Range("F12").Value = 119
$ws.Range("F13").Value = 2498
$ws.Range("F14").Value = 62
$ws.Range("F15").Value = 29
$ws.Range("F17").Value = 47
$ws.Range("F18").Value = 17
$ws.Range("F21").Value = 593
$ws.Range("F22").Value = 177
$ws.Range("F23").Value = 92
$ws.Range("F26").Value = 2095
$ws.Range("F27").Value = 4191
$ws.Range("F31").Value = 1229
$ws.Range("F32").Value = 243
$ws.Range("F33").Value = 2130
$ws.Range("F35").Value = 472
$ws.Range("F37").Value = 126
$ws.Range("F39").Value = 438
$ws.Range("F40").Value = 725
$ws.Range("F43").Value = 7
